$d = $word.ActiveDocument

# Update Ativação date
$d.Content.Find.Execute("Ativação: 01/01/2025", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 15/07/2025", 2) | Out-Null

# --- Objetivos (PT) paragraph (Paragraphs(6)) ---
$rng = $d.Paragraphs(6).Range
$rng.Find.Execute("Fornecer ao aluno uma visão integrada dos sistemas de tratamento de águas residuárias, incluindo conceitos de caracterização dos efluentes e o impacto de lançamento no corpo receptor, fundamentos dos processos e operações de uma estação de tratamento de efluentes e disposição dos resíduos gerados nas estações.", $true, $false, $false, $false, $false, $true, 1, $false, "Poluição aquática; Controle da poluição dos recursos hídricos; Caracterização dos efluentes e o impacto de lançamento no corpo receptor; Tratamento preliminar; Tratamento biológico; Tratamento combinado; Tratamento e disposição do lodo. Visita supervisionada prevista.", 2) | Out-Null
# --- Objetivos (EN) paragraph (Paragraphs(7)) ---
$rng = $d.Paragraphs(7).Range
$rng.Find.Execute("Provide the student with an integrated view of wastewater treatment systems, including concepts of effluent characterization and the impact of release on the receiving body, fundamentals of the processes and operations of an effluent treatment plant and disposal of waste generated at the stations.", $true, $false, $false, $false, $false, $true, 1, $false, "Aquatic pollution; Control of pollution of water resources; Characterization of effluents and the impact of release on the receiving body; Preliminary treatment; Biological treatment; Combined treatment; Sludge treatment and disposal. Supervised visits.", 2) | Out-Null
# --- Docentes list paragraph (Paragraphs(9)) ---
$rng = $d.Paragraphs(9).Range
$rng.Find.Execute("3380737 - Flávio Teixeira da Silva" + [char]11 + "", $true, $false, $false, $false, $false, $true, 1, $false, "Fornecer ao aluno uma visão integrada dos sistemas de tratamento de águas residuárias, incluindo conceitos de caracterização dos efluentes e o impacto de lançamento no corpo receptor, fundamentos dos processos e operações de uma estação de tratamento de efluentes e disposição dos resíduos gerados nas estações." + [char]11 + "", 2) | Out-Null
$rng.Find.Execute("1720367 - Teresa Cristina Brazil de Paiva", $true, $false, $false, $false, $false, $true, 1, $false, "Poluição aquática e controle da poluição; Caracterização física, química e biológica das águas residuais; Métodos de detecção de toxicidade e o impacto do lançamento de efluentes nos corpos receptores; Tratamento preliminar; Princípios da microbiologia do tratamento biológico e ecologia microbiana; Sistemas de lagoas de estabilização; sistemas de lodos ativados; Sistemas com biodisco; Reatores aeróbios e anaeróbicos; Sistemas combinados; Tratamento e disposição final do lodo de estação de tratamento de efluentes. Visita supervisionada a laboratórios e indústrias, a depender da viabilidade no momento do oferecimento da disciplina.", 2) | Out-Null
# --- Programa resumido (PT) paragraph (Paragraphs(11)) ---
$rng = $d.Paragraphs(11).Range
$rng.Find.Execute("Poluição aquática; Controle da poluição dos recursos hídricos; Caracterização dos efluentes e o impacto de lançamento no corpo receptor; Tratamento preliminar; Tratamento biológico; Tratamento combinado; Tratamento e disposição do lodo.", $true, $false, $false, $false, $false, $true, 1, $false, "Aulas expositivas, atividades de projeto e solução de exercícios.", 2) | Out-Null
# --- Programa resumido (EN) paragraph (Paragraphs(12)) ---
$rng = $d.Paragraphs(12).Range
$rng.Find.Execute("Aquatic pollution; Control of pollution of water resources; Characterization of effluents and the impact of release on the receiving body; Preliminary treatment; Biological treatment; Combined treatment; Sludge treatment and disposal.", $true, $false, $false, $false, $false, $true, 1, $false, "Provide the student with an integrated view of wastewater treatment systems, including concepts of effluent characterization and the impact of release on the receiving body, fundamentals of the processes and operations of an effluent treatment plant and disposal of waste generated at the stations.", 2) | Out-Null
# --- Programa (PT) paragraph (Paragraphs(14)) ---
$rng = $d.Paragraphs(14).Range
$rng.Find.Execute("Poluição aquática e controle da poluição; Caracterização física, química e biológica das águas residuais; Métodos de detecção de toxicidade e o impacto do lançamento de efluentes nos corpos receptores; Tratamento preliminar; Princípios da microbiologia do tratamento biológico e ecologia microbiana; Sistemas de lagoas de estabilização; sistemas de lodos ativados; Sistemas com biodisco; Reatores aeróbios e anaeróbicos; Sistemas combinados; Tratamento e disposição final do lodo de estação de tratamento de efluentes.", $true, $false, $false, $false, $false, $true, 1, $false, "Média ponderada das notas atribuídas à prova, exercício, seminário e/ou relatório. Serão aprovados os alunos que obtenham média igual ou maior que 5,0 e 70% de frequência no curso.", 2) | Out-Null
# --- Programa (EN) paragraph (Paragraphs(15)) ---
$rng = $d.Paragraphs(15).Range
$rng.Find.Execute("Aquatic pollution and pollution control; Physical, chemical and biological characterization of wastewater; Methods for detecting toxicity and the impact of releasing effluents into receiving bodies; Preliminary treatment; Principles of microbiology of biological treatment and microbial ecology; Stabilization pond systems; activated sludge systems; Systems with biodisk; Aerobic and anaerobic reactors; Combined systems; Treatment and final disposal of sludge from an effluent treatment plant.", $true, $false, $false, $false, $false, $true, 1, $false, "Aquatic pollution and pollution control; Physical, chemical and biological characterization of wastewater; Methods for detecting toxicity and the impact of releasing effluents into receiving bodies; Preliminary treatment; Principles of microbiology of biological treatment and microbial ecology; Stabilization pond systems; activated sludge systems; Systems with biodisk; Aerobic and anaerobic reactors; Combined systems; Treatment and final disposal of sludge from an effluent treatment plant. Supervised visits to laboratories and industries, depending on feasibility at the time the discipline is offered.", 2) | Out-Null
# --- Avaliação paragraph (Paragraphs(17)) ---
$rng = $d.Paragraphs(17).Range
$rng.Find.Execute("Aulas expositivas, atividades de projeto e solução de exercícios. " + [char]11 + "Viaje didática (visita técnica a ETEs)." + [char]11 + "", $true, $false, $false, $false, $false, $true, 1, $false, "Aos alunos que obtiverem média igual ou maior que 3,0 e menor que 5,0 será oferecido um programa de recuperação que será avaliado por uma prova final. Nesse caso, a média final do aluno será: Média final = (média do período letivo normal + nota prova final)/2. Serão aprovados os alunos que obtiverem média final igual ou maior que 5,0." + [char]11 + "", 2) | Out-Null
$rng.Find.Execute("Média ponderada das notas atribuídas à prova, exercício, seminário e/ou relatório. Serão aprovados os alunos que obtenham média igual ou maior que 5,0 e 70% de frequência no curso." + [char]11 + "", $true, $false, $false, $false, $false, $true, 1, $false, "1.ANDRÉOLI, CV; VON SPERLING, M; FERNANDES, F. Lodo de esgoto: tratamento e disposição final - Princípios do tratamento biológico de águas residuárias. V. 6. - Belo Horizonte: Departamento de Engenharia Sanitária e Ambiental; Universidade Federal de Minas Gerais; 2001." + [char]11 + "2.BITTON, G. Wastewater Microbiology. Willey-Liss (John Wiley and Sons Inc., Publications), 3 ed., 2005." + [char]11 + "3.CHERNICHARO, CAL Reatores anaeróbios - Princípios do tratamento biológico de águas residuárias. V. 5. - Belo Horizonte: Departamento de Engenharia Sanitária e Ambiental; Universidade Federal de Minas Gerais; 2008." + [char]11 + "4.DEZOTTI, M.; SANT'ANNA JUNIOR, G.L.; BASSIN, J.P. Processos Biológicos Avançados para Tratamento de eﬂuentes e técnicas de biologia molecular para o estudo da diversidade microbiana. Rio de Janeiro: Interciência, 2011. 368p." + [char]11 + "5.JORDÃO, E. P.; PESSOA, C. A. Tratamento de esgotos domésticos. 4. ed. Rio de Janeiro: Editora SEGRAC, 2005, 906 pp." + [char]11 + "6.METCALF & EDDY, INC. Wastewater engineering treatment, disposal and reuse. 4th ed. Boston: McGraw-Hill, c2003. 1819 p. (McGraw-Hill series in civil and environmental engineering)." + [char]11 + "7.VON SPERLING, M. Lagoas de estabilização - Princípios do tratamento biológico de águas residuárias. V. 3. 2 ed.- Belo Horizonte: Departamento de Engenharia Sanitária e Ambiental; Universidade Federal de Minas Gerais; 2002." + [char]11 + "8.VON SPERLING, M. Lodos ativados - Princípios do tratamento biológico de águas residuárias. V. 4. 4 ed.- Belo Horizonte: Departamento de Engenharia Sanitária e Ambiental; Universidade Federal de Minas Gerais; 2016." + [char]11 + "9.VON SPERLING, M. Introdução à qualidade das águas e ao tratamento de esgotos - Princípios do tratamento biológico de águas residuárias. V. 1, 4 ed. - Belo Horizonte: Departamento de Engenharia Sanitária e Ambiental; Universidade Federal de Minas Gerais; 2014." + [char]11 + "10.Fugita, R. S. Fundamentos do controle de poluição das águas. CETESB, 2018." + [char]11 + "", 2) | Out-Null
$rng.Find.Execute("Aos alunos que obtiverem média igual ou maior que 3,0 e menor que 5,0 será oferecido um programa de recuperação que será avaliado por uma prova final. Nesse caso, a média final do aluno será: Média final = (média do período letivo normal + nota prova final)/2. Serão aprovados os alunos que obtiverem média final igual ou maior que 5,0.", $true, $false, $false, $false, $false, $true, 1, $false, "3380737 - Flávio Teixeira da Silva", 2) | Out-Null
# --- Bibliografia paragraph (Paragraphs(19)) ---
$rng = $d.Paragraphs(19).Range
$rng.Find.Execute("1.ANDRÉOLI, CV; VON SPERLING, M; FERNANDES, F. Lodo de esgoto: tratamento e disposição final - Princípios do tratamento biológico de águas residuárias. V. 6. - Belo Horizonte: Departamento de Engenharia Sanitária e Ambiental; Universidade Federal de Minas Gerais; 2001." + [char]11 + "2.BITTON, G. Wastewater Microbiology. Willey-Liss (John Wiley and Sons Inc., Publications), 3 ed., 2005." + [char]11 + "3.CHERNICHARO, CAL Reatores anaeróbios - Princípios do tratamento biológico de águas residuárias. V. 5. - Belo Horizonte: Departamento de Engenharia Sanitária e Ambiental; Universidade Federal de Minas Gerais; 2008." + [char]11 + "4.DEZOTTI, M.; SANT'ANNA JUNIOR, G.L.; BASSIN, J.P. Processos Biológicos Avançados para Tratamento de eﬂuentes e técnicas de biologia molecular para o estudo da diversidade microbiana. Rio de Janeiro: Interciência, 2011. 368p." + [char]11 + "5.JORDÃO, E. P.; PESSOA, C. A. Tratamento de esgotos domésticos. 4. ed. Rio de Janeiro: Editora SEGRAC, 2005, 906 pp." + [char]11 + "6.METCALF & EDDY, INC. Wastewater engineering treatment, disposal and reuse. 4th ed. Boston: McGraw-Hill, c2003. 1819 p. (McGraw-Hill series in civil and environmental engineering)." + [char]11 + "7.VON SPERLING, M. Lagoas de estabilização - Princípios do tratamento biológico de águas residuárias. V. 3. 2 ed.- Belo Horizonte: Departamento de Engenharia Sanitária e Ambiental; Universidade Federal de Minas Gerais; 2002." + [char]11 + "8.VON SPERLING, M. Lodos ativados - Princípios do tratamento biológico de águas residuárias. V. 4. 4 ed.- Belo Horizonte: Departamento de Engenharia Sanitária e Ambiental; Universidade Federal de Minas Gerais; 2016." + [char]11 + "9.VON SPERLING, M. Introdução à qualidade das águas e ao tratamento de esgotos - Princípios do tratamento biológico de águas residuárias. V. 1, 4 ed. - Belo Horizonte: Departamento de Engenharia Sanitária e Ambiental; Universidade Federal de Minas Gerais; 2014." + [char]11 + "10.Fugita, R. S. Fundamentos do controle de poluição das águas. CETESB, 2018.", $true, $false, $false, $false, $false, $true, 1, $false, "1720367 - Teresa Cristina Brazil de Paiva", 2) | Out-Null

Write-Host "Done"